$d = $word.ActiveDocument

# Locate the paragraph that ends with "...create it then save to it."
# and insert two brand-new paragraphs right after it (and before the
# pre-existing blank paragraph that follows it).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*create it then save to it.*") {
        $target = $p
    }
}

# --- First new paragraph -------------------------------------------------
$target.Range.InsertParagraphAfter()
$target = $target.Next()

$xml1 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">One line of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ctypes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> sets the desktop wallpaper when I give it an absolute path, now I just need to get all the information to dynamically set the image I just downloaded as the desktop.</w:t></w:r></w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void]$target.Range.InsertXML($xml1)

# --- Second new paragraph -------------------------------------------------
$target.Range.InsertParagraphAfter()
$target = $target.Next()

$xml2 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Got the path by feeding in the original wallpaper URL and getting the image name this way. Not a fan of this method because it couples the methods together a little more than I’d like but it will do for now.</w:t></w:r></w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void]$target.Range.InsertXML($xml2)

Write-Output "done"
